$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.196.48"
$ws.Range("E2").Value = "  -0.28%  "
$ws.Range("D3").Value = "'1.909.10"
$ws.Range("E3").Value = "  -0.33%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "'326.13"
$ws.Range("E5").Value = "  -0.48%  "
$ws.Range("D6").Value = "'1.002"
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("D7").Value = "'0.4625"
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "'0.3896"
$ws.Range("E8").Value = "  -1.28%  "
$ws.Range("D9").Value = "'0.07884"
$ws.Range("E9").Value = "  -0.72%  "
$ws.Range("D10").Value = "'0.9924"
$ws.Range("E10").Value = "  -0.94%  "
$ws.Range("E11").Value = "  -1.70%  "
$ws.Range("D12").Value = "'1.835.86"
$ws.Range("E12").Value = "  -4.72%  "
$ws.Range("D13").Value = "'5.770"
$ws.Range("E13").Value = "  -0.02%  "
$ws.Range("D14").Value = "'7.055"
$ws.Range("E14").Value = "  -0.74%  "
$ws.Range("D15").Value = "'0.07033"
$ws.Range("E15").Value = "  +1.12%  "
$ws.Range("D16").Value = "'88.22"
$ws.Range("E16").Value = "  -0.46%  "
$ws.Range("E17").Value = "  -0.08%  "
$ws.Range("D18").Value = "'0.000009958"
$ws.Range("E18").Value = "  -1.45%  "
$ws.Range("E19").Value = "  -0.16%  "
$ws.Range("E20").Value = "  -0.01%  "
$ws.Range("D21").Value = "'29.214.29"
$ws.Range("E21").Value = "  -0.23%  "
$ws.Range("D22").Value = "'5.332"
$ws.Range("E22").Value = "  -0.53%  "
$ws.Range("E23").Value = "  +0.33%  "
$ws.Range("D24").Value = "'2.128.06"
$ws.Range("E24").Value = "  -1.00%  "
$ws.Range("D25").Value = "'2.096"
$ws.Range("E25").Value = "  +1.46%  "
$ws.Range("D26").Value = "'156.48"
$ws.Range("E26").Value = "  -0.33%  "
$ws.Range("D27").Value = "'19.49"
$ws.Range("D28").Value = "'5.918"
$ws.Range("E28").Value = "  -3.18%  "
$ws.Range("D29").Value = "'118.91"
$ws.Range("E29").Value = "  -0.09%  "
$ws.Range("D30").Value = "'1.883"
$ws.Range("E30").Value = "  -5.72%  "
$ws.Range("D31").Value = "'0.09362"
$ws.Range("D32").Value = "'0.8996"
$ws.Range("E32").Value = "  -2.92%  "
$ws.Range("D33").Value = "'5.238"
$ws.Range("E33").Value = "  -2.29%  "
$ws.Range("D34").Value = "'1.326"
$ws.Range("E34").Value = "  -2.31%  "
$ws.Range("D35").Value = "'3.149"
$ws.Range("E35").Value = "  -3.93%  "
$ws.Range("D36").Value = "'0.05806"
$ws.Range("E36").Value = "  -0.62%  "
$ws.Range("D37").Value = "'1.173"
$ws.Range("E37").Value = "  -2.57%  "
$ws.Range("D38").Value = "'0.02093"
$ws.Range("E38").Value = "  -0.92%  "
$ws.Range("E39").Value = "  -0.13%  "
$ws.Range("D40").Value = "'0.5715"
$ws.Range("E40").Value = "  -0.68%  "
$ws.Range("D41").Value = "'7.690"
$ws.Range("E41").Value = "  -3.65%  "
$ws.Range("D42").Value = "'0.1815"
$ws.Range("E42").Value = "  +0.59%  "
$ws.Range("D43").Value = "'9.756"
$ws.Range("E43").Value = "  -2.31%  "
$ws.Range("D44").Value = "'11.88"
$ws.Range("E44").Value = "  -1.08%  "
$ws.Range("D45").Value = "'0.5366"
$ws.Range("E45").Value = "  -1.19%  "
$ws.Range("D46").Value = "'2.178"
$ws.Range("E46").Value = "  -5.41%  "
$ws.Range("E47").Value = "  -0.87%  "
$ws.Range("D48").Value = "'1.846"
$ws.Range("E48").Value = "  -1.90%  "
$ws.Range("D49").Value = "'2.553"
$ws.Range("E49").Value = "  -0.74%  "
$ws.Range("D50").Value = "'113.50"
$ws.Range("E50").Value = "  +0.02%  "
$ws.Range("D51").Value = "'0.2993"
$ws.Range("E51").Value = "  +1.32%  "
